$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 570.5030999408009
$ws.Range("R2").Value = 5134.527899467208
$ws.Range("S2").Value = 0.02070125718402391
$ws.Range("T2").Value = 0.02070125718402392
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 968.5183961020467
$ws.Range("R3").Value = 8716.665564918421
$ws.Range("S3").Value = 0.03514362745311513
$ws.Range("T3").Value = 0.03514362745311513
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 416.6072567106144
$ws.Range("R4").Value = 3749.46531039553
$ws.Range("S4").Value = 0.01511699755319824
$ws.Range("T4").Value = 0.01511699755319825
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 7291.227203713299
$ws.Range("R5").Value = 65621.04483341968
$ws.Range("S5").Value = 0.264569236428133
$ws.Range("T5").Value = 0.264569236428133
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 12378.00053617374
$ws.Range("R6").Value = 111402.0048255637
$ws.Range("S6").Value = 0.4491477303977974
$ws.Range("T6").Value = 0.4491477303977975
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 5324.385027369703
$ws.Range("R7").Value = 47919.46524632732
$ws.Range("S7").Value = 0.1932004643091049
$ws.Range("T7").Value = 0.193200464309105
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 177.8408341374649
$ws.Range("R8").Value = 1600.567507237184
$ws.Range("S8").Value = 0.006453126802786908
$ws.Range("T8").Value = 0.006453126802786912
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 301.9126792792898
$ws.Range("R9").Value = 2717.214113513608
$ws.Range("S9").Value = 0.01095519379566359
$ws.Range("T9").Value = 0.01095519379566359
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 129.8674486586044
$ws.Range("R10").Value = 1168.80703792744
$ws.Range("S10").Value = 0.004712366076176904
$ws.Range("T10").Value = 0.004712366076176906
